$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$cell = $ws.Range('D2')
$cell.NumberFormat = "@"
$cell.Value = '58.165.28'
$cell.Style = "Normal"
$ws.Range('E2').Value = '  +1.80%  '

# Row 3
$cell = $ws.Range('D3')
$cell.NumberFormat = "@"
$cell.Value = '2.364.92'
$cell.Style = "Normal"
$ws.Range('E3').Value = '  +2.05%  '

# Row 4
$ws.Range('E4').Value = '  -0.27%  '

# Row 5
$cell = $ws.Range('D5')
$cell.NumberFormat = "@"
$cell.Value = '544.20'
$cell.Style = "Normal"
$ws.Range('E5').Value = '  +2.62%  '

# Row 6
$cell = $ws.Range('D6')
$cell.NumberFormat = "@"
$cell.Value = '136.38'
$cell.Style = "Normal"
$ws.Range('E6').Value = '  +2.96%  '

# Row 7
$ws.Range('E7').Value = '  +0.48%  '

# Row 8
$ws.Range('E8').Value = '  +5.23%  '

# Row 9
$ws.Range('E9').Value = '  +1.56%  '

# Row 10
$ws.Range('E10').Value = '  +4.00%  '

# Row 11
$ws.Range('E11').Value = '  -0.85%  '

# Row 12
$cell = $ws.Range('D12')
$cell.NumberFormat = "@"
$cell.Value = '0.357'
$cell.Style = "Normal"
$ws.Range('E12').Value = '  +1.51%  '

# Row 13
$cell = $ws.Range('D13')
$cell.NumberFormat = "@"
$cell.Value = '24.02'
$cell.Style = "Normal"
$ws.Range('E13').Value = '  +3.07%  '

# Row 14
$cell = $ws.Range('D14')
$cell.NumberFormat = "@"
$cell.Value = '2.782.57'
$cell.Style = "Normal"
$ws.Range('E14').Value = '  +1.73%  '

# Row 15
$cell = $ws.Range('D15')
$cell.NumberFormat = "@"
$cell.Value = '58.101.33'
$cell.Style = "Normal"
$ws.Range('E15').Value = '  +1.75%  '

# Row 16
$ws.Range('E16').Value = '  +1.95%  '

# Row 17
$cell = $ws.Range('D17')
$cell.NumberFormat = "@"
$cell.Value = '2.378.79'
$cell.Style = "Normal"
$ws.Range('E17').Value = '  +2.02%  '

# Row 18
$ws.Range('E18').Value = '  +3.63%  '

# Row 19
$cell = $ws.Range('D19')
$cell.NumberFormat = "@"
$cell.Value = '333.40'
$cell.Style = "Normal"
$ws.Range('E19').Value = '  -1.06%  '

# Row 20
$ws.Range('E20').Value = '  +2.58%  '

# Row 21
$cell = $ws.Range('D21')
$cell.NumberFormat = "@"
$cell.Value = '6.79'
$cell.Style = "Normal"
$ws.Range('E21').Value = '  +0.47%  '

# Row 22
$ws.Range('E22').Value = '  +0.28%  '

# Row 23
$cell = $ws.Range('D23')
$cell.NumberFormat = "@"
$cell.Value = '62.90'
$cell.Style = "Normal"
$ws.Range('E23').Value = '  +1.71%  '

# Row 24
$ws.Range('E24').Value = '  +0.42%  '

# Row 25
$cell = $ws.Range('D25')
$cell.NumberFormat = "@"
$cell.Value = '8.53'
$cell.Style = "Normal"
$ws.Range('E25').Value = '  -1.85%  '

# Row 26
$cell = $ws.Range('D26')
$cell.NumberFormat = "@"
$cell.Value = '0.999'
$cell.Style = "Normal"
$ws.Range('E26').Value = '  +0.55%  '

# Row 27
$ws.Range('E27').Value = '  +2.42%  '

# Row 28
$ws.Range('B28').Value = 'PancakeSwap'
$ws.Range('C28').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$cell = $ws.Range('D28')
$cell.NumberFormat = "@"
$cell.Value = '1.76'
$cell.Style = "Normal"
$ws.Range('E28').Value = '  +2.38%  '

# Row 29
$ws.Range('B29').Value = 'Monero'
$ws.Range('C29').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$cell = $ws.Range('D29')
$cell.NumberFormat = "@"
$cell.Value = '172.78'
$cell.Style = "Normal"
$ws.Range('E29').Value = '  -0.31%  '

# Row 30
$ws.Range('E30').Value = '  +2.40%  '

# Row 31
$cell = $ws.Range('D31')
$cell.NumberFormat = "@"
$cell.Value = '6.17'
$cell.Style = "Normal"
$ws.Range('E31').Value = '  +1.16%  '

# Row 32
$ws.Range('E32').Value = '  +11.71%  '

# Row 33
$cell = $ws.Range('D33')
$cell.NumberFormat = "@"
$cell.Value = '18.56'
$cell.Style = "Normal"
$ws.Range('E33').Value = '  +0.42%  '

# Row 34
$ws.Range('E34').Value = '  +0.05%  '

# Row 35
$ws.Range('E35').Value = '  +6.83%  '

# Row 36
$cell = $ws.Range('D36')
$cell.NumberFormat = "@"
$cell.Value = '0.999'
$cell.Style = "Normal"
$ws.Range('E36').Value = '  +0.53%  '

# Row 37
$ws.Range('E37').Value = '  +0.58%  '

# Row 38
$ws.Range('E38').Value = '  +3.88%  '

# Row 39
$cell = $ws.Range('D39')
$cell.NumberFormat = "@"
$cell.Value = '39.45'
$cell.Style = "Normal"
$ws.Range('E39').Value = '  +0.64%  '

# Row 40
$cell = $ws.Range('D40')
$cell.NumberFormat = "@"
$cell.Value = '145.53'
$cell.Style = "Normal"
$ws.Range('E40').Value = '  -2.61%  '

# Row 41
$cell = $ws.Range('D41')
$cell.NumberFormat = "@"
$cell.Value = '293.81'
$cell.Style = "Normal"
$ws.Range('E41').Value = '  +3.43%  '

# Row 42
$cell = $ws.Range('D42')
$cell.NumberFormat = "@"
$cell.Value = '0.379'
$cell.Style = "Normal"
$ws.Range('E42').Value = '  +1.22%  '

# Row 43
$ws.Range('E43').Value = '  +1.43%  '

# Row 44
$cell = $ws.Range('D44')
$cell.NumberFormat = "@"
$cell.Value = '0.0948'
$cell.Style = "Normal"
$ws.Range('E44').Value = '  +1.98%  '

# Row 45
$ws.Range('E45').Value = '  +3.06%  '

# Row 46
$ws.Range('E46').Value = '  +1.06%  '

# Row 47
$cell = $ws.Range('D47')
$cell.NumberFormat = "@"
$cell.Value = '0.566'
$cell.Style = "Normal"
$ws.Range('E47').Value = '  +1.24%  '

# Row 48
$ws.Range('E48').Value = '  +3.01%  '

# Row 49
$cell = $ws.Range('D49')
$cell.NumberFormat = "@"
$cell.Value = '17.55'
$cell.Style = "Normal"
$ws.Range('E49').Value = '  +0.81%  '

# Row 50
$ws.Range('E50').Value = '  +0.11%  '

# Row 51
$cell = $ws.Range('D51')
$cell.NumberFormat = "@"
$cell.Value = '11.07'
$cell.Style = "Normal"
$ws.Range('E51').Value = '  +0.47%  '

Write-Host "applied cryptos update"